$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text cells keep their literal string values (avoid numeric/date auto-conversion)
$cellData = @{
    'D2' = '68.181.14'
    'E2' = '  -0.50%  '
    'D3' = '3.910.22'
    'E3' = '  -1.41%  '
    'E4' = '  -0.19%  '
    'D5' = '483.48'
    'E5' = '  -0.01%  '
    'D6' = '146.38'
    'E6' = '  -3.34%  '
    'D7' = '0.621'
    'E7' = '  -0.85%  '
    'E8' = '  +0.00%  '
    'D9' = '0.732'
    'E9' = '  -0.25%  '
    'E10' = '  -0.88%  '
    'E11' = '  -2.19%  '
    'D12' = '43.09'
    'E12' = '  -1.31%  '
    'D13' = '10.68'
    'E13' = '  +1.95%  '
    'D14' = '4.537.62'
    'E14' = '  -1.52%  '
    'D15' = '3.933.23'
    'E15' = '  -0.76%  '
    'D16' = '14.20'
    'E16' = '  -3.99%  '
    'E17' = '  -0.77%  '
    'D18' = '20.17'
    'E18' = '  +0.69%  '
    'D19' = '1.14'
    'E19' = '  -0.41%  '
    'D20' = '68.209.62'
    'E20' = '  -0.58%  '
    'D21' = '429.27'
    'E21' = '  -1.99%  '
    'D22' = '3.50'
    'E22' = '  +4.02%  '
    'D23' = '15.14'
    'E23' = '  +5.17%  '
    'D24' = '88.30'
    'E24' = '  +0.13%  '
    'D25' = '11.69'
    'E25' = '  +18.86%  '
    'D26' = '3.69'
    'E26' = '  +1.93%  '
    'E27' = '  +10.63%  '
    'D28' = '37.70'
    'E28' = '  -2.00%  '
    'D29' = '5.66'
    'E29' = '  -1.82%  '
    'D30' = '716.67'
    'E30' = '  -1.56%  '
    'D31' = '13.82'
    'E31' = '  +4.38%  '
    'D32' = '0.131'
    'E32' = '  +2.52%  '
    'D33' = '2.92'
    'E33' = '  +3.22%  '
    'D34' = '0.0₃0905'
    'E34' = '  +5.49%  '
    'D35' = '6.21'
    'E35' = '  +15.33%  '
    'D36' = '41.56'
    'E36' = '  -2.92%  '
    'D37' = '60.45'
    'E37' = '  +0.70%  '
    'B38' = 'Dai'
    'C38' = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
    'D38' = '0.999'
    'E38' = '  +0.11%  '
    'E39' = '  +18.02%  '
    'B40' = 'Fetch.AI'
    'C40' = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
    'D40' = '3.00'
    'E40' = '  +15.87%  '
    'B41' = 'Kaspa'
    'C41' = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
    'D41' = '0.142'
    'E41' = '  -5.52%  '
    'E42' = '  +3.16%  '
    'E43' = '  +2.89%  '
    'D44' = '3.01'
    'E44' = '  +3.86%  '
    'E45' = '  +0.63%  '
    'D46' = '3.33'
    'E46' = '  +2.57%  '
    'E47' = '  -0.09%  '
    'D48' = '3.44'
    'E48' = '  +0.84%  '
    'E49' = '  -3.52%  '
    'D50' = '144.67'
    'E50' = '  -2.52%  '
    'E51' = '  +27.46%  '
}

foreach ($ref in $cellData.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $cellData[$ref]
}
